# ============================================================
# Edit: Update latest output (run 127)
# Rewrites "Schedule" rows 2-5 and "Detailed" rows 11-97
# per the optimisation re-run described in the commit diff.
# ============================================================
$wb = $excel.ActiveWorkbook

# ---------- Sheet "Schedule": rows 2-5 ----------
$wsSchedule = $wb.Worksheets.Item("Schedule")

$scheduleData = @(
  @{R=2; A=46043; B=46043.1875; C=4.5; D=17.01; E=573.8929364999999; F=33.73856181657848},
  @{R=3; A=46043.29166666666; B=46043.66666666666; C=9; D=34.02; E=-221.8333455; F=-6.520674470899471},
  @{R=4; A=46043.89583333334; B=46044.10416666666; C=5; D=18.9; E=588.914547; F=31.15949984126984},
  @{R=5; A=46044.27083333334; B=46044.66666666666; C=9.5; D=35.91; E=13.66849575; F=0.3806320175438597}
)

foreach ($row in $scheduleData) {
    $r = $row.R
    $wsSchedule.Cells.Item($r,1).Value = $row.A
    $wsSchedule.Cells.Item($r,2).Value = $row.B
    $wsSchedule.Cells.Item($r,3).Value = $row.C
    $wsSchedule.Cells.Item($r,4).Value = $row.D
    $wsSchedule.Cells.Item($r,5).Value = $row.E
    $wsSchedule.Cells.Item($r,6).Value = $row.F
}

# Newly added rows 3:5 need the same date/time display format as row 2 (A:B)
$wsSchedule.Range("A3:B5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------- Sheet "Detailed" ----------
$wsDetailed = $wb.Worksheets.Item("Detailed")

# Rows 11-15: Pump_Status flips ON -> OFF
foreach ($r in @(11,12,13,14,15)) {
    $wsDetailed.Cells.Item($r,5).Value = "OFF"
}

# Rows 21-49: updated Price (B), Type (C) and Pump_Status (E) from the new optimisation run
$detailedUpdates = @(
  @{R=21; B=-7.41772; C="historical"; E="ON"},
  @{R=22; B=-9.5; C="historical"; E="ON"},
  @{R=23; B=-10; C="historical"; E="ON"},
  @{R=24; B=-14; C="historical"; E="ON"},
  @{R=25; B=-15.89865; C="historical"; E="ON"},
  @{R=26; B=-22.10072; C="historical"; E="ON"},
  @{R=27; B=-23.5; C="historical"; E="ON"},
  @{R=28; B=-24.41017; C="historical"; E="ON"},
  @{R=29; B=-23.5; C="historical"; E="ON"},
  @{R=30; B=-27; C="historical"; E="ON"},
  @{R=31; B=-25.94511; C="historical"; E="ON"},
  @{R=32; B=-24.08764; C="historical"; E="ON"},
  @{R=33; B=-22.86107; C="historical"; E="ON"},
  @{R=34; B=-6.8; C="historical"; E="OFF"},
  @{R=35; B=-5.51; C="historical"; E="OFF"},
  @{R=36; B=36.06; C="historical"; E="OFF"},
  @{R=37; B=36.08141; C="historical"; E="OFF"},
  @{R=38; B=46.54611; C="historical"; E="OFF"},
  @{R=39; B=64.35683; C="forecast"; E="OFF"},
  @{R=40; B=73.19; C="forecast"; E="OFF"},
  @{R=41; B=77.94; C="forecast"; E="OFF"},
  @{R=42; B=100.01; C="forecast"; E="OFF"},
  @{R=43; B=74.43841999999999; C="forecast"; E="OFF"},
  @{R=44; B=73.19; C="forecast"; E="OFF"},
  @{R=45; B=71.7281; C="forecast"; E="ON"},
  @{R=46; B=59.62291; C="forecast"; E="ON"},
  @{R=47; B=57.98348; C="forecast"; E="ON"},
  @{R=48; B=60.73801; C="forecast"; E="ON"},
  @{R=49; B=61.32156; C="forecast"; E="ON"}
)

foreach ($row in $detailedUpdates) {
    $r = $row.R
    $wsDetailed.Cells.Item($r,2).Value = $row.B
    $wsDetailed.Cells.Item($r,3).Value = $row.C
    $wsDetailed.Cells.Item($r,5).Value = $row.E
}

# Rows 50-97: brand-new rows appended for the next day (2026-01-22) of the forecast horizon
$detailedNewRows = @(
  @{R=50; A=46044; B=57.31; C="forecast"; D=46044; E="ON"},
  @{R=51; A=46044.02083333334; B=63.88086; C="forecast"; D=46044; E="ON"},
  @{R=52; A=46044.04166666666; B=57.31; C="forecast"; D=46044; E="ON"},
  @{R=53; A=46044.0625; B=57.06; C="forecast"; D=46044; E="ON"},
  @{R=54; A=46044.08333333334; B=57.06; C="forecast"; D=46044; E="ON"},
  @{R=55; A=46044.10416666666; B=57.06; C="forecast"; D=46044; E="OFF"},
  @{R=56; A=46044.125; B=63.73519; C="forecast"; D=46044; E="OFF"},
  @{R=57; A=46044.14583333334; B=64.10364; C="forecast"; D=46044; E="OFF"},
  @{R=58; A=46044.16666666666; B=65.32088; C="forecast"; D=46044; E="OFF"},
  @{R=59; A=46044.1875; B=65.85026999999999; C="forecast"; D=46044; E="OFF"},
  @{R=60; A=46044.20833333334; B=66.04559; C="forecast"; D=46044; E="OFF"},
  @{R=61; A=46044.22916666666; B=76.12006; C="forecast"; D=46044; E="OFF"},
  @{R=62; A=46044.25; B=64.89; C="forecast"; D=46044; E="OFF"},
  @{R=63; A=46044.27083333334; B=57.06; C="forecast"; D=46044; E="ON"},
  @{R=64; A=46044.29166666666; B=30.9379; C="forecast"; D=46044; E="ON"},
  @{R=65; A=46044.3125; B=0.51; C="forecast"; D=46044; E="ON"},
  @{R=66; A=46044.33333333334; B=-5.50985; C="forecast"; D=46044; E="ON"},
  @{R=67; A=46044.35416666666; B=-6.15086; C="forecast"; D=46044; E="ON"},
  @{R=68; A=46044.375; B=-10; C="forecast"; D=46044; E="ON"},
  @{R=69; A=46044.39583333334; B=-13.50737; C="forecast"; D=46044; E="ON"},
  @{R=70; A=46044.41666666666; B=-9.621499999999999; C="forecast"; D=46044; E="ON"},
  @{R=71; A=46044.4375; B=-12.01; C="forecast"; D=46044; E="ON"},
  @{R=72; A=46044.45833333334; B=-12.01; C="forecast"; D=46044; E="ON"},
  @{R=73; A=46044.47916666666; B=-5.74313; C="forecast"; D=46044; E="ON"},
  @{R=74; A=46044.5; B=-10; C="forecast"; D=46044; E="ON"},
  @{R=75; A=46044.52083333334; B=-8.0564; C="forecast"; D=46044; E="ON"},
  @{R=76; A=46044.54166666666; B=-7.79393; C="forecast"; D=46044; E="ON"},
  @{R=77; A=46044.5625; B=-5.88864; C="forecast"; D=46044; E="ON"},
  @{R=78; A=46044.58333333334; B=-5.27725; C="forecast"; D=46044; E="ON"},
  @{R=79; A=46044.60416666666; B=0.51; C="forecast"; D=46044; E="ON"},
  @{R=80; A=46044.625; B=0.51; C="forecast"; D=46044; E="ON"},
  @{R=81; A=46044.64583333334; B=36.06; C="forecast"; D=46044; E="ON"},
  @{R=82; A=46044.66666666666; B=0.51; C="forecast"; D=46044; E="OFF"},
  @{R=83; A=46044.6875; B=-4.13512; C="forecast"; D=46044; E="OFF"},
  @{R=84; A=46044.70833333334; B=-5.14805; C="forecast"; D=46044; E="OFF"},
  @{R=85; A=46044.72916666666; B=-6.90848; C="forecast"; D=46044; E="OFF"},
  @{R=86; A=46044.75; B=-0.45834; C="forecast"; D=46044; E="OFF"},
  @{R=87; A=46044.77083333334; B=0.00036; C="forecast"; D=46044; E="OFF"},
  @{R=88; A=46044.79166666666; B=10.48193; C="forecast"; D=46044; E="OFF"},
  @{R=89; A=46044.8125; B=55.33036; C="forecast"; D=46044; E="OFF"},
  @{R=90; A=46044.83333333334; B=53.90468; C="forecast"; D=46044; E="OFF"},
  @{R=91; A=46044.85416666666; B=54.47327; C="forecast"; D=46044; E="OFF"},
  @{R=92; A=46044.875; B=57.01318; C="forecast"; D=46044; E="OFF"},
  @{R=93; A=46044.89583333334; B=50.38252; C="forecast"; D=46044; E="OFF"},
  @{R=94; A=46044.91666666666; B=30.67112; C="forecast"; D=46044; E="OFF"},
  @{R=95; A=46044.9375; B=56.98; C="forecast"; D=46044; E="OFF"},
  @{R=96; A=46044.95833333334; B=56.23018; C="forecast"; D=46044; E="OFF"},
  @{R=97; A=46044.97916666666; B=48.31676; C="forecast"; D=46044; E="OFF"}
)

foreach ($row in $detailedNewRows) {
    $r = $row.R
    $wsDetailed.Cells.Item($r,1).Value = $row.A
    $wsDetailed.Cells.Item($r,2).Value = $row.B
    $wsDetailed.Cells.Item($r,3).Value = $row.C
    $wsDetailed.Cells.Item($r,4).Value = $row.D
    $wsDetailed.Cells.Item($r,5).Value = $row.E
}

# New rows need the same DateTime / Date display formats as the existing data (A:D)
$wsDetailed.Range("A50:A97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsDetailed.Range("D50:D97").NumberFormat = "YYYY-MM-DD"
